# The only substantive content change in the target revision is a
# one-letter typo fix in the next-to-last paragraph:
#   "...on the client in out local storage..."
# becomes
#   "...on the client in our local storage..."
#
# (Everything else in the recorded diff is Word's automatic
# spell-check/grammar-check "proofing" markup -- <w:proofErr/> tags and
# the purely cosmetic run-splitting that comes along with it -- which
# Word regenerates on the fly from its live proofing engine rather than
# from any deliberate editor action, so we only need to reproduce the
# actual wording change here.)

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "client in out local storage",  # old text (FindText)
    $true,                          # MatchCase
    $false,                         # MatchWholeWord
    $false,                         # MatchWildcards
    $false,                         # MatchSoundsLike
    $false,                         # MatchAllWordForms
    $true,                          # Forward
    1,                              # Wrap (wdFindContinue)
    $false,                         # Format
    "client in our local storage",  # ReplaceWith
    2                               # Replace (wdReplaceAll)
)

Write-Output "replaced: $found"
